$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "GenEntity": mark the first 4 detail rows (rows 3-6) of the first
# table as Primary Key (column I) -- matches the other already-flagged rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("GenEntity")
$ws1.Range("I3").Value = 1
$ws1.Range("I4").Value = 1
$ws1.Range("I5").Value = 1
$ws1.Range("I6").Value = 1

# ---------------------------------------------------------------------------
# Sheet "DataTypeMapping": insert 3 new rows right under the header for the
# "null handling" mapping (null / Null / NULL -> String), which pushes the
# existing From/To/Library rows down by 3, then add a new small auto "Data
# Type -> Prefix" lookup table (+ truncated-number count) in columns E:H.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DataTypeMapping")

$ws2.Rows("2:4").Insert()

$ws2.Range("A2").Value = "null"
$ws2.Range("B2").Value = "String"
$ws2.Range("A3").Value = "Null"
$ws2.Range("B3").Value = "String"
$ws2.Range("A4").Value = "NULL"
$ws2.Range("B4").Value = "String"
$ws2.Range("A2:B4").Font.Bold = $false

# New "Data Type" / "Prefix" helper table used to auto fill the Prefix
# column, plus a running count of truncated-number conversions.
$ws2.Range("E1").Value = "Data Type"
$ws2.Range("F1").Value = "Prefix"
$ws2.Range("H1").Value = "Truncated Num"
$ws2.Range("E1:F1").Font.Bold = $true
$ws2.Range("H1").Font.Bold = $true

$ws2.Range("E2").Value = "String"
$ws2.Range("F2").Value = "str"
$ws2.Range("H2").Value = 3

$ws2.Range("E3").Value = "int"
$ws2.Range("F3").Value = "int"

$ws2.Range("E4").Value = "Integer"
$ws2.Range("F4").Value = "int"

$ws2.Range("E5").Value = "double"
$ws2.Range("F5").Value = "dbl"

$ws2.Range("E6").Value = "Double"
$ws2.Range("F6").Value = "dbl"

$ws2.Range("E7").Value = "float"
$ws2.Range("F7").Value = "flt"

$ws2.Range("E8").Value = "Float"
$ws2.Range("F8").Value = "flt"

$ws2.Range("E9").Value = "short"
$ws2.Range("F9").Value = "sht"

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A12").Select()

$ws2.Activate()
$ws2.Range("E16").Select()
